$wb = $excel.ActiveWorkbook

# ---- Sheet 1: LP1912 ----
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: 11:33:52"
$ws1.Range("A3").Value = "Total filas: 139"

# Rows 66-67 swap
$ws1.Cells.Item(66, 1).Value = "08:11:18"
$ws1.Cells.Item(66, 2).Value = "09:28"
$ws1.Cells.Item(66, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(66, 4).Value = 77
$ws1.Cells.Item(66, 5).Value = "LP1912"
$ws1.Cells.Item(67, 1).Value = "08:28:52"
$ws1.Cells.Item(67, 2).Value = "09:28"
$ws1.Cells.Item(67, 3).Value = "10_OLMOS"
$ws1.Cells.Item(67, 4).Value = 60
$ws1.Cells.Item(67, 5).Value = "LP1912"

# Rows 117-144 rewritten (includes 10 new rows 135-144)
$ws1.Cells.Item(117, 1).Value = "11:33:52"
$ws1.Cells.Item(117, 2).Value = "11:52"
$ws1.Cells.Item(117, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(117, 4).Value = 19
$ws1.Cells.Item(117, 5).Value = "LP1912"
$ws1.Cells.Item(118, 1).Value = "10:36:50"
$ws1.Cells.Item(118, 2).Value = "11:52"
$ws1.Cells.Item(118, 3).Value = "225_GOMEZ"
$ws1.Cells.Item(118, 4).Value = 76
$ws1.Cells.Item(118, 5).Value = "LP1912"
$ws1.Cells.Item(119, 1).Value = "10:04:30"
$ws1.Cells.Item(119, 2).Value = "11:53"
$ws1.Cells.Item(119, 3).Value = "225_GOMEZ"
$ws1.Cells.Item(119, 4).Value = 109
$ws1.Cells.Item(119, 5).Value = "LP1912"
$ws1.Cells.Item(120, 1).Value = "10:04:30"
$ws1.Cells.Item(120, 2).Value = "11:58"
$ws1.Cells.Item(120, 3).Value = "17_ROMERO"
$ws1.Cells.Item(120, 4).Value = 114
$ws1.Cells.Item(120, 5).Value = "LP1912"
$ws1.Cells.Item(121, 1).Value = "10:36:50"
$ws1.Cells.Item(121, 2).Value = "12:05"
$ws1.Cells.Item(121, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(121, 4).Value = 89
$ws1.Cells.Item(121, 5).Value = "LP1912"
$ws1.Cells.Item(122, 1).Value = "10:56:15"
$ws1.Cells.Item(122, 2).Value = "12:06"
$ws1.Cells.Item(122, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(122, 4).Value = 70
$ws1.Cells.Item(122, 5).Value = "LP1912"
$ws1.Cells.Item(123, 1).Value = "10:36:50"
$ws1.Cells.Item(123, 2).Value = "12:10"
$ws1.Cells.Item(123, 3).Value = "15_ABASTO"
$ws1.Cells.Item(123, 4).Value = 94
$ws1.Cells.Item(123, 5).Value = "LP1912"
$ws1.Cells.Item(124, 1).Value = "10:36:50"
$ws1.Cells.Item(124, 2).Value = "12:10"
$ws1.Cells.Item(124, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(124, 4).Value = 94
$ws1.Cells.Item(124, 5).Value = "LP1912"
$ws1.Cells.Item(125, 1).Value = "11:33:52"
$ws1.Cells.Item(125, 2).Value = "12:16"
$ws1.Cells.Item(125, 3).Value = "10_OLMOS"
$ws1.Cells.Item(125, 4).Value = 43
$ws1.Cells.Item(125, 5).Value = "LP1912"
$ws1.Cells.Item(126, 1).Value = "11:13:15"
$ws1.Cells.Item(126, 2).Value = "12:17"
$ws1.Cells.Item(126, 3).Value = "10_OLMOS"
$ws1.Cells.Item(126, 4).Value = 64
$ws1.Cells.Item(126, 5).Value = "LP1912"
$ws1.Cells.Item(127, 1).Value = "10:36:50"
$ws1.Cells.Item(127, 2).Value = "12:21"
$ws1.Cells.Item(127, 3).Value = "215C_EL PATO"
$ws1.Cells.Item(127, 4).Value = 105
$ws1.Cells.Item(127, 5).Value = "LP1912"
$ws1.Cells.Item(128, 1).Value = "10:56:15"
$ws1.Cells.Item(128, 2).Value = "12:22"
$ws1.Cells.Item(128, 3).Value = "215C_EL PATO"
$ws1.Cells.Item(128, 4).Value = 86
$ws1.Cells.Item(128, 5).Value = "LP1912"
$ws1.Cells.Item(129, 1).Value = "11:13:15"
$ws1.Cells.Item(129, 2).Value = "12:29"
$ws1.Cells.Item(129, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(129, 4).Value = 76
$ws1.Cells.Item(129, 5).Value = "LP1912"
$ws1.Cells.Item(130, 1).Value = "10:36:50"
$ws1.Cells.Item(130, 2).Value = "12:32"
$ws1.Cells.Item(130, 3).Value = "14_ABASTO"
$ws1.Cells.Item(130, 4).Value = 116
$ws1.Cells.Item(130, 5).Value = "LP1912"
$ws1.Cells.Item(131, 1).Value = "11:33:52"
$ws1.Cells.Item(131, 2).Value = "12:32"
$ws1.Cells.Item(131, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(131, 4).Value = 59
$ws1.Cells.Item(131, 5).Value = "LP1912"
$ws1.Cells.Item(132, 1).Value = "10:56:15"
$ws1.Cells.Item(132, 2).Value = "12:33"
$ws1.Cells.Item(132, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(132, 4).Value = 97
$ws1.Cells.Item(132, 5).Value = "LP1912"
$ws1.Cells.Item(133, 1).Value = "10:56:15"
$ws1.Cells.Item(133, 2).Value = "12:33"
$ws1.Cells.Item(133, 3).Value = "14_ABASTO"
$ws1.Cells.Item(133, 4).Value = 97
$ws1.Cells.Item(133, 5).Value = "LP1912"
$ws1.Cells.Item(134, 1).Value = "10:36:50"
$ws1.Cells.Item(134, 2).Value = "12:34"
$ws1.Cells.Item(134, 3).Value = "15_ABASTO"
$ws1.Cells.Item(134, 4).Value = 118
$ws1.Cells.Item(134, 5).Value = "LP1912"
$ws1.Cells.Item(135, 1).Value = "10:49:38"
$ws1.Cells.Item(135, 2).Value = "12:36"
$ws1.Cells.Item(135, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(135, 4).Value = 107
$ws1.Cells.Item(135, 5).Value = "LP1912"
$ws1.Cells.Item(136, 1).Value = "11:33:52"
$ws1.Cells.Item(136, 2).Value = "12:47"
$ws1.Cells.Item(136, 3).Value = "14_ABASTO"
$ws1.Cells.Item(136, 4).Value = 74
$ws1.Cells.Item(136, 5).Value = "LP1912"
$ws1.Cells.Item(137, 1).Value = "10:49:38"
$ws1.Cells.Item(137, 2).Value = "12:48"
$ws1.Cells.Item(137, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(137, 4).Value = 119
$ws1.Cells.Item(137, 5).Value = "LP1912"
$ws1.Cells.Item(138, 1).Value = "11:33:52"
$ws1.Cells.Item(138, 2).Value = "12:48"
$ws1.Cells.Item(138, 3).Value = "15X38_ABASTO"
$ws1.Cells.Item(138, 4).Value = 75
$ws1.Cells.Item(138, 5).Value = "LP1912"
$ws1.Cells.Item(139, 1).Value = "11:33:52"
$ws1.Cells.Item(139, 2).Value = "13:02"
$ws1.Cells.Item(139, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(139, 4).Value = 89
$ws1.Cells.Item(139, 5).Value = "LP1912"
$ws1.Cells.Item(140, 1).Value = "11:13:15"
$ws1.Cells.Item(140, 2).Value = "13:03"
$ws1.Cells.Item(140, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(140, 4).Value = 110
$ws1.Cells.Item(140, 5).Value = "LP1912"
$ws1.Cells.Item(141, 1).Value = "11:33:52"
$ws1.Cells.Item(141, 2).Value = "13:03"
$ws1.Cells.Item(141, 3).Value = "215C_EL PATO"
$ws1.Cells.Item(141, 4).Value = 90
$ws1.Cells.Item(141, 5).Value = "LP1912"
$ws1.Cells.Item(142, 1).Value = "11:33:52"
$ws1.Cells.Item(142, 2).Value = "13:13"
$ws1.Cells.Item(142, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(142, 4).Value = 100
$ws1.Cells.Item(142, 5).Value = "LP1912"
$ws1.Cells.Item(143, 1).Value = "11:33:52"
$ws1.Cells.Item(143, 2).Value = "13:17"
$ws1.Cells.Item(143, 3).Value = "10_OLMOS"
$ws1.Cells.Item(143, 4).Value = 104
$ws1.Cells.Item(143, 5).Value = "LP1912"
$ws1.Cells.Item(144, 1).Value = "11:33:52"
$ws1.Cells.Item(144, 2).Value = "13:25"
$ws1.Cells.Item(144, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(144, 4).Value = 112
$ws1.Cells.Item(144, 5).Value = "LP1912"

# ---- Sheet 2: LP1912-215 ----
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 11:33:52"
$ws2.Range("A3").Value = "Total filas: 25"
$ws2.Cells.Item(30, 1).Value = "11:33:52"
$ws2.Cells.Item(30, 2).Value = "13:03"
$ws2.Cells.Item(30, 3).Value = "215C_EL PATO"
$ws2.Cells.Item(30, 4).Value = 90
$ws2.Cells.Item(30, 5).Value = "LP1912"

# ---- Sheet 3: 6203-6173 ----
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 11:33:52"
$ws3.Range("A3").Value = "Total filas: 22"
$ws3.Cells.Item(25, 1).Value = "11:33:52"
$ws3.Cells.Item(25, 2).Value = "13:11"
$ws3.Cells.Item(25, 3).Value = "215C_LA PLATA"
$ws3.Cells.Item(25, 4).Value = 98
$ws3.Cells.Item(25, 5).Value = "L6203"
$ws3.Cells.Item(26, 1).Value = "11:13:15"
$ws3.Cells.Item(26, 2).Value = "13:12"
$ws3.Cells.Item(26, 3).Value = "215C_LA PLATA"
$ws3.Cells.Item(26, 4).Value = 119
$ws3.Cells.Item(26, 5).Value = "L6203"
$ws3.Cells.Item(27, 1).Value = "11:33:52"
$ws3.Cells.Item(27, 2).Value = "13:20"
$ws3.Cells.Item(27, 3).Value = "215B_LP-P MOR-1 Y 57"
$ws3.Cells.Item(27, 4).Value = 107
$ws3.Cells.Item(27, 5).Value = "L6173"
